$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 changes
$ws.Range("Q3").Value = 1.84
$ws.Range("R3").Value = 1.89

# Row 4 changes
$ws.Range("G4").Value = 1.65
$ws.Range("H4").Value = 3.8
$ws.Range("I4").Value = 5.25
$ws.Range("J4").Value = 2.3
$ws.Range("L4").Value = 6
$ws.Range("U4").Value = 3.3
$ws.Range("V4").Value = 1.32
$ws.Range("AA4").Value = 2.2
$ws.Range("AB4").Value = 1.62
$ws.Range("AC4").Value = 5.5
$ws.Range("AJ4").Value = 7.5
$ws.Range("AL4").Value = 81
$ws.Range("AO4").Value = 26

# Row 5 changes
$ws.Range("G5").Value = 4.65
$ws.Range("H5").Value = 3.25
$ws.Range("I5").Value = 1.75
$ws.Range("J5").Value = 4.75
$ws.Range("K5").Value = 2.02
$ws.Range("L5").Value = 2.4
$ws.Range("O5").Value = 1.28
$ws.Range("P5").Value = 3.05
$ws.Range("S5").Value = 1.82
$ws.Range("T5").Value = 1.78
$ws.Range("W5").Value = 2.87
$ws.Range("X5").Value = 1.31
$ws.Range("Z5").Value = 2.47
$ws.Range("AA5").Value = 1.72
$ws.Range("AB5").Value = 1.9
$ws.Range("AC5").Value = 14
$ws.Range("AD5").Value = 30
$ws.Range("AE5").Value = 14.5
$ws.Range("AF5").Value = 90
$ws.Range("AH5").Value = 40
$ws.Range("AI5").Value = 9.5
$ws.Range("AJ5").Value = 6.4
$ws.Range("AK5").Value = 13.5
$ws.Range("AL5").Value = 60
$ws.Range("AM5").Value = 450
$ws.Range("AN5").Value = 6.7
$ws.Range("AO5").Value = 8.25
$ws.Range("AP5").Value = 8
$ws.Range("AQ5").Value = 14.5
$ws.Range("AR5").Value = 14.5
$ws.Range("AS5").Value = 25
